# Append a new traffic-data row (row 9) to the sheet, matching the
# existing rows: Location, Date (stored as text, not an Excel date
# serial), cars-up, cars-down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value = "Location 1"

# Force the date column to be stored as literal text "2024-12-13"
# (matching B2:B8) instead of letting Excel auto-convert it to a date
# serial number. Resetting the style back to Normal afterwards keeps
# the cell's style index at the sheet default (no explicit s="...")
# just like the other data rows.
$ws.Range("B9").NumberFormat = "@"
$ws.Range("B9").Value = "2024-12-13"
$ws.Range("B9").Style = "Normal"

$ws.Range("C9").Value = 24
$ws.Range("D9").Value = 34
